$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width: H gets a bit wider (bestFit no longer matches, now custom) ---
$ws.Columns.Item(8).ColumnWidth = 12.67

# --- Re-align the whole header + data row block (A1:Z2) to left+vcenter ---
# (keeps each column's existing font / wrap / number-format, just adds
# horizontal=left on top of them, and also "fills in" any previously
# unused cells in the rectangle with the same style)
$ws.Range("A1:Z2").HorizontalAlignment = -4131

# --- Drop the leftover QUERY/EXPL_QUERY headers (R1:W1) and the
#     USER_DB/PASSWORD_DB/HOSTNAME headers (X1:Z1) ---
$ws.Range("R1:Z1").ClearContents()

# --- Drop the raw SQL query cell ---
$ws.Range("R2").ClearContents()

# --- Drop the helper formula cell ---
$ws.Range("U2").ClearContents()

# --- Drop the raw DB credentials (sa / password / host) ---
$ws.Range("X2:Z2").ClearContents()

# --- Restore plain A1 selection / scroll position (drop stale
#     topLeftCell=K1 / selection=N2 leftover view state) ---
$ws.Range("A1").Select()
